$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new value looks like a plain number (single decimal point)
# must be forced to Text format first, otherwise Excel auto-converts the literal
# string into a floating point number (losing exact text / introducing fp noise).
$textCells = @("D5", "D6", "D10", "D12", "D14", "D21", "D23", "D24", "D26", "D27", "D28", "D29", "D30", "D31", "D32", "D33", "D35", "D36", "D37", "D38", "D39", "D41", "D42", "D45", "D47", "D48", "D50", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "69.437.93"
$ws.Range("E2").Value = "  -0.97%  "
$ws.Range("D3").Value = "3.514.15"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").Value = "574.49"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").Value = "185.79"
$ws.Range("E6").Value = "  -2.87%  "
$ws.Range("D7").Value = "3.503.23"
$ws.Range("E7").Value = "  -2.46%  "
$ws.Range("E8").Value = "  -2.92%  "
$ws.Range("E9").Value = "  +0.04%  "
$ws.Range("D10").Value = "0.189"
$ws.Range("E10").Value = "  +4.40%  "
$ws.Range("E11").Value = "  -2.28%  "
$ws.Range("D12").Value = "54.24"
$ws.Range("E12").Value = "  -3.14%  "
$ws.Range("E13").Value = "  -1.33%  "
$ws.Range("D14").Value = "9.48"
$ws.Range("E14").Value = "  -1.99%  "
$ws.Range("D15").Value = "4.078.70"
$ws.Range("E15").Value = "  -2.23%  "
$ws.Range("E16").Value = "  -2.60%  "
$ws.Range("D17").Value = "69.330.60"
$ws.Range("E17").Value = "  -1.10%  "
$ws.Range("D18").Value = "3.518.00"
$ws.Range("E18").Value = "  -1.93%  "
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").Value = "  -1.04%  "
$ws.Range("D21").Value = "547.54"
$ws.Range("E21").Value = "  +14.93%  "
$ws.Range("E22").Value = "  -2.92%  "
$ws.Range("D23").Value = "18.73"
$ws.Range("E23").Value = "  -2.51%  "
$ws.Range("D24").Value = "4.97"
$ws.Range("E24").Value = "  -1.21%  "
$ws.Range("E25").Value = "  +1.61%  "
$ws.Range("D26").Value = "94.39"
$ws.Range("E26").Value = "  -0.92%  "
$ws.Range("D27").Value = "11.44"
$ws.Range("E27").Value = "  +3.55%  "
$ws.Range("D28").Value = "2.95"
$ws.Range("E28").Value = "  -1.90%  "
$ws.Range("D29").Value = "9.17"
$ws.Range("E29").Value = "  -1.62%  "
$ws.Range("D30").Value = "31.94"
$ws.Range("E30").Value = "  -1.13%  "
$ws.Range("D31").Value = "7.27"
$ws.Range("E31").Value = "  -5.45%  "
$ws.Range("D32").Value = "12.66"
$ws.Range("E32").Value = "  +3.49%  "
$ws.Range("D33").Value = "64.84"
$ws.Range("E34").Value = "  -4.26%  "
$ws.Range("D35").Value = "547.23"
$ws.Range("E35").Value = "  -6.70%  "
$ws.Range("D36").Value = "0.406"
$ws.Range("E36").Value = "  +2.48%  "
$ws.Range("B37").Value = "Fetch.AI"
$ws.Range("C37").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D37").Value = "3.09"
$ws.Range("E37").Value = "  +8.32%  "
$ws.Range("B38").Value = "InjectiveProtocol"
$ws.Range("C38").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D38").Value = "38.14"
$ws.Range("E38").Value = "  -2.63%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "0.999"
$ws.Range("E39").Value = "  -0.27%  "
$ws.Range("E40").Value = "  -4.34%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "3.11"
$ws.Range("E41").Value = "  -3.63%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "3.38"
$ws.Range("E42").Value = "  -2.41%  "
$ws.Range("E43").Value = "  -3.23%  "
$ws.Range("D44").Value = "3.300.86"
$ws.Range("E44").Value = "  +2.27%  "
$ws.Range("D45").Value = "3.00"
$ws.Range("E45").Value = "  -2.36%  "
$ws.Range("E46").Value = "  +0.26%  "
$ws.Range("D47").Value = "3.45"
$ws.Range("E47").Value = "  +3.07%  "
$ws.Range("D48").Value = "0.134"
$ws.Range("E48").Value = "  -2.55%  "
$ws.Range("E49").Value = "  -5.44%  "
$ws.Range("D50").Value = "0.999"
$ws.Range("E50").Value = "  +0.06%  "
$ws.Range("D51").Value = "137.09"
$ws.Range("E51").Value = "  +2.29%  "

# Restore the default style on the cells we temporarily reformatted, so only the
# cell content changes and no stray style index is introduced.
foreach ($addr in $textCells) {
    $ws.Range($addr).Style = "Normal"
}
